{"js": "const replacements = [\n  [\"2025-11-11 Tuesday\", \"2025-11-12 Wednesday\"],\n  [\"537\u00d76=\", \"329\u00d77=\"],\n  [\"459\u00d72=\", \"725\u00d75=\"],\n  [\"629\u00d76=\", \"678\u00d73=\"],\n  [\"593\u00d72=\", \"310\u00d78=\"],\n  [\"218\u00d72=\", \"419\u00d72=\"],\n  [\"640\u00d73=\", \"661\u00d72=\"],\n  [\"510\u00d76=\", \"301\u00d74=\"],\n  [\"972\u00d75=\", \"695\u00d75=\"],\n  [\"801\u00d72=\", \"979\u00d78=\"],\n  [\"746\u00d78=\", \"745\u00d79=\"],\n  [\"291\u00d78=\", \"872\u00d73=\"],\n  [\"773\u00d76=\", \"393\u00d72=\"],\n  [\"110\u00d75=\", \"115\u00d75=\"],\n  [\"772\u00d79=\", \"593\u00d74=\"],\n  [\"484\u00d77=\", \"901\u00d79=\"],\n  [\"982\u00d76=\", \"431\u00d77=\"],\n  [\"505\u00d72=\", \"571\u00d74=\"],\n  [\"498\u00d72=\", \"194\u00d78=\"],\n  [\"305\u00d74=\", \"816\u00d73=\"],\n  [\"346\u00d75=\", \"596\u00d76=\"],\n  [\"752\u00d77=\", \"804\u00d72=\"],\n  [\"588\u00d78=\", \"694\u00d73=\"],\n  [\"110\u00d72=\", \"988\u00d75=\"],\n  [\"707\u00d78=\", \"987\u00d73=\"],\n  [\"275\u00d76=\", \"644\u00d75=\"],\n];\n\nconst body = context.document.body;\nfor (const [find, replace] of replacements) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  for (const r of results.items) {\n    r.insertText(replace, 'Replace');\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2025-11-11 Tuesday\", \"2025-11-12 Wednesday\"),\n  @(\"537\u00d76=\", \"329\u00d77=\"),\n  @(\"459\u00d72=\", \"725\u00d75=\"),\n  @(\"629\u00d76=\", \"678\u00d73=\"),\n  @(\"593\u00d72=\", \"310\u00d78=\"),\n  @(\"218\u00d72=\", \"419\u00d72=\"),\n  @(\"640\u00d73=\", \"661\u00d72=\"),\n  @(\"510\u00d76=\", \"301\u00d74=\"),\n  @(\"972\u00d75=\", \"695\u00d75=\"),\n  @(\"801\u00d72=\", \"979\u00d78=\"),\n  @(\"746\u00d78=\", \"745\u00d79=\"),\n  @(\"291\u00d78=\", \"872\u00d73=\"),\n  @(\"773\u00d76=\", \"393\u00d72=\"),\n  @(\"110\u00d75=\", \"115\u00d75=\"),\n  @(\"772\u00d79=\", \"593\u00d74=\"),\n  @(\"484\u00d77=\", \"901\u00d79=\"),\n  @(\"982\u00d76=\", \"431\u00d77=\"),\n  @(\"505\u00d72=\", \"571\u00d74=\"),\n  @(\"498\u00d72=\", \"194\u00d78=\"),\n  @(\"305\u00d74=\", \"816\u00d73=\"),\n  @(\"346\u00d75=\", \"596\u00d76=\"),\n  @(\"752\u00d77=\", \"804\u00d72=\"),\n  @(\"588\u00d78=\", \"694\u00d73=\"),\n  @(\"110\u00d72=\", \"988\u00d75=\"),\n  @(\"707\u00d78=\", \"987\u00d73=\"),\n  @(\"275\u00d76=\", \"644\u00d75=\"),\n)\n\nforeach ($pair in $pairs) {\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Text = $pair[0]\n  $find.Replacement.ClearFormatting()\n  $find.Replacement.Text = $pair[1]\n  $find.Execute(\n    $pair[0],   # FindText\n    $true,      # MatchCase\n    $false,     # MatchWholeWord\n    $false,     # MatchWildcards\n    $false,     # MatchSoundsLike\n    $false,     # MatchAllWordForms\n    $true,      # Forward\n    1,          # Wrap (wdFindContinue)\n    $false,     # Format\n    $pair[1],   # ReplaceWith\n    2           # Replace (wdReplaceAll)\n  )\n}\n"}
